$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Value" column (C) had three transcription errors (misplaced/missing
# digits) that threw off the "Total check" row (C17 = SUM(C2:C15)-C16,
# which should net to 0 when the port totals agree with the line items).
$ws.Range("C2").Value = 376183
$ws.Range("C3").Value = 83620
$ws.Range("C16").Value = 687383

# Recalculate so C17/D17's cached formula results reflect the corrected
# inputs (C17 goes from -75258 to 0).
$excel.Calculate()

# Leave the sheet scrolled/selected where the author ended up after
# verifying the fix: window scrolled to row 24, cell C17 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C17").Select()
